$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11 : "image trop lourdes" issue
$ws.Range("B11").Value = "image trop lourdes"
$ws.Range("C11").Value = "utilisation d'un format non adapté"
$ws.Range("D11").Value = "utiliser le format adapter"
$ws.Range("E11").Value = "changer les formats aux formats adaptés"

# New row 12 : "fichier non minifier" issue
$ws.Range("B12").Value = "fichier non minifier"
$ws.Range("C12").Value = "fichier plus lourd qu'ils ne le devraient"
$ws.Range("D12").Value = "minifier les fichiers une fois terminés"

# E12 repeats the "Action recommandée" text from D12 (same as several rows above),
# so copy D12 into E12 to carry over both value and formatting.
$ws.Range("D12").Copy($ws.Range("E12"))

# Columns D and E grew wider (best-fit) to accommodate the new, longer text.
$ws.Columns.Item(4).ColumnWidth = 29.77734375
$ws.Columns.Item(5).ColumnWidth = 33.33203125

# Update the active selection left behind after the edit.
$ws.Range("E12").Select() | Out-Null
